$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# New header cells D1/E1 for the two new JML checklist columns
$ws.Range("D1").Value = "LOOPINVARIANT"
$ws.Range("E1").Value = "PRIVATE INVARIANT"

# Match the bold+underline header style used by the other header cells (A1:C1, F1)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("D1:E1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New data cells on the MinMaxStrategy row, same "KLAAR" marker as its neighbours
$ws.Range("B17").Value = "KLAAR"
$ws.Range("C17").Value = "KLAAR"

# Size the two new columns to fit their header text (as Excel's own AutoFit would)
$ws.Columns.Item(4).ColumnWidth = 15.166666666666666
$ws.Columns.Item(5).ColumnWidth = 18.25

# Selection moves to the row that was just edited
$ws.Range("C17").Select() | Out-Null
